$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44159
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value = 7000
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = 3000
$ws.Range("O2").Value = 'Provincia de Chacabuco'
$ws.Range("P2").Value = 30

# Row 3
$ws.Range("D3").Value = 44845
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value = 7900
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = 3000
$ws.Range("O3").Value = 'Provincia de Chacabuco'
$ws.Range("P3").Value = 30

# Row 4
$ws.Range("D4").Value = 44215
$ws.Range("I4").Value = 'Primera'
$ws.Range("J4").Value = 16000
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 3000
$ws.Range("O4").Value = 'Provincia de Chacabuco'
$ws.Range("P4").Value = 30

# Row 5
$ws.Range("D5").Value = 44189
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 16000
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = 3000
$ws.Range("O5").Value = 'Provincia de Chacabuco'
$ws.Range("P5").Value = 30

# Row 6
$ws.Range("D6").Value = 44167
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 7000
$ws.Range("K6").Value = 3000
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = 3000
$ws.Range("O6").Value = 'Provincia de Chacabuco'
$ws.Range("P6").Value = 30

# Row 7
$ws.Range("D7").Value = 44855
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 7900
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = 3000
$ws.Range("O7").Value = 'Provincia de Chacabuco'
$ws.Range("P7").Value = 30

# Row 8
$ws.Range("D8").Value = 44214
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 7000
$ws.Range("K8").Value = 3000
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = 3000
$ws.Range("O8").Value = 'Provincia de Chacabuco'
$ws.Range("P8").Value = 30

# Row 9
$ws.Range("D9").Value = 44188
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 12000
$ws.Range("K9").Value = 3000
$ws.Range("L9").Value = 3000
$ws.Range("M9").Value = 3000
$ws.Range("O9").Value = 'Provincia de Chacabuco'
$ws.Range("P9").Value = 30

# Row 10
$ws.Range("D10").Value = 44209
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 7000
$ws.Range("K10").Value = 2500
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 2750
$ws.Range("O10").Value = 'Provincia de Chacabuco'
$ws.Range("P10").Value = 28

# Row 11
$ws.Range("D11").Value = 44902
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 7000
$ws.Range("K11").Value = 3000
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = 3000
$ws.Range("O11").Value = 'Provincia de Chacabuco'
$ws.Range("P11").Value = 30

# Row 12
$ws.Range("D12").Value = 44230
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 16000
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("M12").Value = 3000
$ws.Range("O12").Value = 'Provincia de Chacabuco'
$ws.Range("P12").Value = 30

# Row 13
$ws.Range("D13").Value = 44602
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 12000
$ws.Range("K13").Value = 3000
$ws.Range("L13").Value = 3000
$ws.Range("M13").Value = 3000
$ws.Range("O13").Value = 'Provincia de Chacabuco'
$ws.Range("P13").Value = 30

# Row 14
$ws.Range("D14").Value = 44602
$ws.Range("I14").Value = 'Segunda'
$ws.Range("J14").Value = 6000
$ws.Range("K14").Value = 2500
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2500
$ws.Range("O14").Value = 'Provincia de Chacabuco'
$ws.Range("P14").Value = 25

# Row 15
$ws.Range("D15").Value = 44210
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value = 8800
$ws.Range("K15").Value = 2500
$ws.Range("L15").Value = 3000
$ws.Range("M15").Value = 2750
$ws.Range("O15").Value = 'Provincia de Chacabuco'
$ws.Range("P15").Value = 28

# Row 16
$ws.Range("D16").Value = 44161
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 7000
$ws.Range("K16").Value = 3000
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = 3000
$ws.Range("O16").Value = 'Provincia de Chacabuco'
$ws.Range("P16").Value = 30

# Row 17
$ws.Range("D17").Value = 44245
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 9000
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = 3000
$ws.Range("O17").Value = 'Región Metropolitana'
$ws.Range("P17").Value = 30

# Row 18
$ws.Range("D18").Value = 44245
$ws.Range("I18").Value = 'Segunda'
$ws.Range("J18").Value = 5000
$ws.Range("K18").Value = 2500
$ws.Range("L18").Value = 2500
$ws.Range("M18").Value = 2500
$ws.Range("O18").Value = 'Región Metropolitana'
$ws.Range("P18").Value = 25

# Row 19
$ws.Range("D19").Value = 44186
$ws.Range("I19").Value = 'Primera'
$ws.Range("J19").Value = 10000
$ws.Range("K19").Value = 3000
$ws.Range("L19").Value = 3000
$ws.Range("M19").Value = 3000
$ws.Range("O19").Value = 'Provincia de Chacabuco'
$ws.Range("P19").Value = 30

# Row 20
$ws.Range("D20").Value = 44873
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 7900
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 3000
$ws.Range("M20").Value = 3000
$ws.Range("O20").Value = 'Provincia de Chacabuco'
$ws.Range("P20").Value = 30

# Row 21
$ws.Range("D21").Value = 44882
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value = 7900
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = 3000
$ws.Range("M21").Value = 3000
$ws.Range("O21").Value = 'Provincia de Chacabuco'
$ws.Range("P21").Value = 30

# Row 22
$ws.Range("D22").Value = 44881
$ws.Range("I22").Value = 'Primera'
$ws.Range("J22").Value = 7900
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = 3000
$ws.Range("O22").Value = 'Provincia de Chacabuco'
$ws.Range("P22").Value = 30

# Row 23
$ws.Range("D23").Value = 44229
$ws.Range("I23").Value = 'Primera'
$ws.Range("J23").Value = 16000
$ws.Range("K23").Value = 3000
$ws.Range("L23").Value = 3000
$ws.Range("M23").Value = 3000
$ws.Range("O23").Value = 'Provincia de Chacabuco'
$ws.Range("P23").Value = 30

# Row 24
$ws.Range("D24").Value = 44880
$ws.Range("I24").Value = 'Primera'
$ws.Range("J24").Value = 7900
$ws.Range("K24").Value = 3000
$ws.Range("L24").Value = 3000
$ws.Range("M24").Value = 3000
$ws.Range("O24").Value = 'Provincia de Chacabuco'
$ws.Range("P24").Value = 30

# Row 25
$ws.Range("D25").Value = 44876
$ws.Range("I25").Value = 'Primera'
$ws.Range("J25").Value = 7900
$ws.Range("K25").Value = 3000
$ws.Range("L25").Value = 3000
$ws.Range("M25").Value = 3000
$ws.Range("O25").Value = 'Provincia de Chacabuco'
$ws.Range("P25").Value = 30

# Row 26
$ws.Range("D26").Value = 44166
$ws.Range("I26").Value = 'Primera'
$ws.Range("J26").Value = 7000
$ws.Range("K26").Value = 3000
$ws.Range("L26").Value = 3000
$ws.Range("M26").Value = 3000
$ws.Range("O26").Value = 'Provincia de Chacabuco'
$ws.Range("P26").Value = 30

# Row 27
$ws.Range("D27").Value = 44847
$ws.Range("I27").Value = 'Primera'
$ws.Range("J27").Value = 7900
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = 3000
$ws.Range("O27").Value = 'Provincia de Chacabuco'
$ws.Range("P27").Value = 30

# Row 28
$ws.Range("D28").Value = 44204
$ws.Range("I28").Value = 'Primera'
$ws.Range("J28").Value = 7000
$ws.Range("K28").Value = 3000
$ws.Range("L28").Value = 3000
$ws.Range("M28").Value = 3000
$ws.Range("O28").Value = 'Provincia de Chacabuco'
$ws.Range("P28").Value = 30

# Row 29
$ws.Range("D29").Value = 44859
$ws.Range("I29").Value = 'Primera'
$ws.Range("J29").Value = 7900
$ws.Range("K29").Value = 3000
$ws.Range("L29").Value = 3000
$ws.Range("M29").Value = 3000
$ws.Range("O29").Value = 'Provincia de Chacabuco'
$ws.Range("P29").Value = 30

# Row 30
$ws.Range("D30").Value = 44875
$ws.Range("I30").Value = 'Primera'
$ws.Range("J30").Value = 7900
$ws.Range("K30").Value = 3000
$ws.Range("L30").Value = 3000
$ws.Range("M30").Value = 3000
$ws.Range("O30").Value = 'Provincia de Chacabuco'
$ws.Range("P30").Value = 30

# Row 31
$ws.Range("D31").Value = 44168
$ws.Range("I31").Value = 'Primera'
$ws.Range("J31").Value = 7000
$ws.Range("K31").Value = 3000
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = 3000
$ws.Range("O31").Value = 'Provincia de Chacabuco'
$ws.Range("P31").Value = 30

# Row 32
$ws.Range("D32").Value = 44860
$ws.Range("I32").Value = 'Primera'
$ws.Range("J32").Value = 7900
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = 3000
$ws.Range("O32").Value = 'Provincia de Chacabuco'
$ws.Range("P32").Value = 30

# Row 33
$ws.Range("D33").Value = 44874
$ws.Range("I33").Value = 'Primera'
$ws.Range("J33").Value = 7900
$ws.Range("K33").Value = 3000
$ws.Range("L33").Value = 3000
$ws.Range("M33").Value = 3000
$ws.Range("O33").Value = 'Provincia de Chacabuco'
$ws.Range("P33").Value = 30

# Row 34
$ws.Range("D34").Value = 44162
$ws.Range("I34").Value = 'Primera'
$ws.Range("J34").Value = 7000
$ws.Range("K34").Value = 3000
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = 3000
$ws.Range("O34").Value = 'Provincia de Chacabuco'
$ws.Range("P34").Value = 30

# Row 35
$ws.Range("D35").Value = 44231
$ws.Range("I35").Value = 'Primera'
$ws.Range("J35").Value = 12000
$ws.Range("K35").Value = 3000
$ws.Range("L35").Value = 3000
$ws.Range("M35").Value = 3000
$ws.Range("O35").Value = 'Provincia de Chacabuco'
$ws.Range("P35").Value = 30

# Row 36
$ws.Range("D36").Value = 44160
$ws.Range("I36").Value = 'Primera'
$ws.Range("J36").Value = 7000
$ws.Range("K36").Value = 3000
$ws.Range("L36").Value = 3000
$ws.Range("M36").Value = 3000
$ws.Range("O36").Value = 'Provincia de Chacabuco'
$ws.Range("P36").Value = 30

# Row 37
$ws.Range("D37").Value = 44846
$ws.Range("I37").Value = 'Primera'
$ws.Range("J37").Value = 7900
$ws.Range("K37").Value = 3000
$ws.Range("L37").Value = 3000
$ws.Range("M37").Value = 3000
$ws.Range("O37").Value = 'Provincia de Chacabuco'
$ws.Range("P37").Value = 30

# Row 38
$ws.Range("D38").Value = 44901
$ws.Range("I38").Value = 'Primera'
$ws.Range("J38").Value = 7000
$ws.Range("K38").Value = 3000
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = 3000
$ws.Range("O38").Value = 'Provincia de Chacabuco'
$ws.Range("P38").Value = 30

# Row 39
$ws.Range("D39").Value = 44232
$ws.Range("I39").Value = 'Primera'
$ws.Range("J39").Value = 16000
$ws.Range("K39").Value = 3000
$ws.Range("L39").Value = 3000
$ws.Range("M39").Value = 3000
$ws.Range("O39").Value = 'Provincia de Chacabuco'
$ws.Range("P39").Value = 30

# Row 40
$ws.Range("D40").Value = 44600
$ws.Range("I40").Value = 'Primera'
$ws.Range("J40").Value = 1300
$ws.Range("K40").Value = 3500
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = 3808
$ws.Range("O40").Value = 'Región Metropolitana'
$ws.Range("P40").Value = 38

# Row 41
$ws.Range("D41").Value = 44181
$ws.Range("I41").Value = 'Primera'
$ws.Range("J41").Value = 12000
$ws.Range("K41").Value = 3000
$ws.Range("L41").Value = 3000
$ws.Range("M41").Value = 3000
$ws.Range("O41").Value = 'Provincia de Chacabuco'
$ws.Range("P41").Value = 30

# Row 42
$ws.Range("D42").Value = 44187
$ws.Range("I42").Value = 'Primera'
$ws.Range("J42").Value = 12000
$ws.Range("K42").Value = 3000
$ws.Range("L42").Value = 3000
$ws.Range("M42").Value = 3000
$ws.Range("O42").Value = 'Provincia de Chacabuco'
$ws.Range("P42").Value = 30

# Row 43
$ws.Range("D43").Value = 44883
$ws.Range("I43").Value = 'Primera'
$ws.Range("J43").Value = 9700
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 3000
$ws.Range("M43").Value = 3000
$ws.Range("O43").Value = 'Provincia de Chacabuco'
$ws.Range("P43").Value = 30

